$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Expand the "7 PSO practices" mention in the Site Lead bullet
# ---------------------------------------------------------------------------
$find1 = $d.Content.Find
$ok1 = $find1.Execute(
    "Serve as Site Lead overseeing delivery governance across all 7 PSO practices in Southeast Asia",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Serve as Site Lead overseeing delivery governance across all 7 PSO practices (Data Analytics, AI/ML, Infrastructure, Security, Enterprise Architecture, Application Development, Delivery Management) in Southeast Asia",
    2)
Write-Host "Replace PSO practices (Site Lead bullet): $ok1"

# ---------------------------------------------------------------------------
# 2) Expand the "7 PSO practices" mention in the agentic AI bullet
# ---------------------------------------------------------------------------
$find2 = $d.Content.Find
$ok2 = $find2.Execute(
    "Pioneered agentic AI adoption across 7 PSO practices and 6 JAPAC sub-regions",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Pioneered agentic AI adoption across all 7 PSO practices (Data Analytics, AI/ML, Infrastructure, Security, Enterprise Architecture, Application Development, Delivery Management) and 6 JAPAC sub-regions",
    2)
Write-Host "Replace PSO practices (agentic AI bullet): $ok2"

# ---------------------------------------------------------------------------
# 3) Rename the "Technical Innovation & Research (Official IP):" heading
# ---------------------------------------------------------------------------
$find3 = $d.Content.Find
$ok3 = $find3.Execute(
    "Technical Innovation & Research (Official IP):",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Technical Innovation & Research:",
    2)
Write-Host "Rename heading: $ok3"

# ---------------------------------------------------------------------------
# 4) Consolidate the six "Technical Innovation & Research" bullets into two
# ---------------------------------------------------------------------------
# Locate the bullet list: find the heading paragraph, then the following
# six ListParagraph bullets that make up the "Official IP" list.
$paras = $d.Paragraphs
$headingIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Text.StartsWith("Technical Innovation & Research:")) {
        $headingIndex = $i
        break
    }
}
Write-Host "Heading paragraph index: $headingIndex"

$bullet1 = $paras.Item($headingIndex + 1)
$bullet2 = $paras.Item($headingIndex + 2)
$bullet3 = $paras.Item($headingIndex + 3)
$bullet4 = $paras.Item($headingIndex + 4)
$bullet5 = $paras.Item($headingIndex + 5)
$bullet6 = $paras.Item($headingIndex + 6)

Write-Host "Bullet1: $($bullet1.Range.Text)"
Write-Host "Bullet6: $($bullet6.Range.Text)"

# Rewrite bullet 1 text (excluding the trailing paragraph mark)
$r1 = $bullet1.Range
$r1b = $d.Range($r1.Start, $r1.End - 1)
$r1b.Text = "5 Google Technical Disclosures on AI and distributed systems - UPIR (automated system synthesis, 274x speedup), FTCS (context architecture for AI agents), ARTEMIS (multi-agent debate framework), ETLC (data processing for GenAI), and LLM inference optimization (speculative decoding, custom Triton kernels)."

# Rewrite bullet 2 text (excluding the trailing paragraph mark)
$r2 = $bullet2.Range
$r2b = $d.Range($r2.Start, $r2.End - 1)
$r2b.Text = "Industry-agnostic agentic AI for enterprise trust decisions. APLS self-learning + cascade routing achieving 86% cost reduction, sub-50ms latency. Won Google Cloud PSO Hackathon JAPAC, qualified for World Finals."

# Delete bullets 3-6 entirely (their paragraph marks included)
$delRange = $d.Range($bullet3.Range.Start, $bullet6.Range.End)
$delRange.Delete()

Write-Host "Done."
